$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price/Volume range as Text so the numeric- and percent-
# looking literals below are stored as strings (matching the source data),
# not auto-converted to Number/Percentage by Excel's type inference.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "300.99"
$ws.Range("E2").Value = "-4.21%"
$ws.Range("D3").Value = "35.55"
$ws.Range("E3").Value = "-1.46%"
$ws.Range("D4").Value = "5.044"
$ws.Range("E4").Value = "-1.39%"
$ws.Range("D5").Value = "0.07988"
$ws.Range("E5").Value = "-1.81%"
$ws.Range("D6").Value = "1.920"
$ws.Range("E6").Value = "-9.53%"
$ws.Range("D7").Value = "7.811"
$ws.Range("E7").Value = "-1.72%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "4.064"
$ws.Range("E8").Value = "-1.80%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9251"
$ws.Range("E9").Value = "-0.72%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1436"
$ws.Range("E10").Value = "38.13%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1898"
$ws.Range("E11").Value = "-1.60%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09247"
$ws.Range("E12").Value = "1.09%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03391"
$ws.Range("E13").Value = "-6.83%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09871"
$ws.Range("E14").Value = "-0.29%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001392"
$ws.Range("E15").Value = "-3.42%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005764"
$ws.Range("E16").Value = "-0.55%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.517"
$ws.Range("E17").Value = "1.44%"
$ws.Range("E18").Value = "5.17%"
$ws.Range("D19").Value = "0.3404"
$ws.Range("E19").Value = "-0.09%"
$ws.Range("E20").Value = "-2.19%"
$ws.Range("D21").Value = "5.055"
$ws.Range("E21").Value = "-1.02%"
$ws.Range("D22").Value = "0.2402"
$ws.Range("E22").Value = "8.54%"
$ws.Range("D23").Value = "0.04490"
$ws.Range("E23").Value = "-1.19%"
$ws.Range("E24").Value = "-2.73%"
$ws.Range("D25").Value = "0.004789"
$ws.Range("E25").Value = "2.32%"
$ws.Range("D26").Value = "0.0001231"
$ws.Range("E26").Value = "-1.81%"
$ws.Range("D27").Value = "0.0003003"
$ws.Range("E27").Value = "-33.44%"
$ws.Range("D39").Value = "0.01915"
$ws.Range("E39").Value = "-2.34%"
$ws.Range("D40").Value = "0.04748"
$ws.Range("E40").Value = "-2.90%"
$ws.Range("D41").Value = "0.007352"
$ws.Range("E41").Value = "-3.27%"
$ws.Range("D42").Value = "0.009658"
$ws.Range("E42").Value = "23.72%"
$ws.Range("D43").Value = "0.1328"
$ws.Range("E43").Value = "-4.20%"
$ws.Range("D44").Value = "0.002111"
$ws.Range("E44").Value = "0.14%"
$ws.Range("D45").Value = "0.01031"
$ws.Range("E45").Value = "-12.25%"
$ws.Range("D46").Value = "0.00006255"
$ws.Range("E46").Value = "-7.41%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.18%"
$ws.Range("E48").Value = "-58.17%"
$ws.Range("D49").Value = "0.001660"
$ws.Range("E49").Value = "-2.58%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.18%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.18%"

# Restore the default cell style so formatting matches the original file
# (only the text content changed, not the styling).
$ws.Range("D2:E51").Style = "Normal"
